$wb = $excel.ActiveWorkbook
$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "4292-MS-EI-DB-SAR-REC-RNI-FEE+INT-FFConMONTHLYonDAY25-FIFC-1-FFROP-DL-FIFR-1-MD-TR-1st"

# Update the product name text wherever it is referenced (input + output
# sheets share the same string) so the trailing "-ONT-PER" becomes "st".
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Update the shortname (B2) from the numeric 4292 to the text "429s",
# removing the test-case inter-dependency on the numeric value.
$wsInput.Range("B2").Value = "429s"

# Move the active selection to B3, matching the post-edit state.
$wsInput.Range("B3").Select()
